# Applies the "automatic update" commit:
#  1. Column C ("Förändrad") date stamp moves from 45184 to 45186 for every
#     data row (rows 2..246).
#  2. Every HYPERLINK() formula in columns S:Y gains a second argument –
#     the friendly display text – equal to the row's designation in
#     column A (e.g. HYPERLINK("...A 37062-2021.xlsx") becomes
#     HYPERLINK("...A 37062-2021.xlsx", "A 37062-2021")).
#     Only rows 2..8 currently contain any HYPERLINK formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (used range) so the script keeps working
# even if the sheet grows/shrinks.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column

# --- 1) Bump the "Förändrad" timestamp (column C) for every data row ------
$oldStamp = 45184
$newStamp = 45186

$changedRange = $ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item($lastRow, 3))
foreach ($cell in $changedRange) {
    if ($cell.Value2 -eq $oldStamp) {
        $cell.Value2 = $newStamp
    }
}

# --- 2) Add the display-text argument to every HYPERLINK() formula --------
for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($label)) { continue }

    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $f = $cell.Formula
        if ([string]::IsNullOrEmpty($f)) { continue }
        if ($f -notlike "*HYPERLINK(*") { continue }
        if ($f -like "*,*") { continue }  # already has a second argument

        $trimmed = $f.TrimEnd()
        if (-not $trimmed.EndsWith(")")) { continue }

        $newFormula = $trimmed.Substring(0, $trimmed.Length - 1) + ', "' + $label + '")'
        $cell.Formula = $newFormula
    }
}
